# daily auto push: 2026-01-21 13:53 UTC
# A new scraped data record (2026/01/21, 水, 19, 27) is inserted at row 671,
# pushing the existing rows 671:712 down to 672:713.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 671, shifting rows 671:712 down to 672:713.
$ws.Rows("671:671").Insert()

# Write the new row's values. Column A holds plain text dates (e.g. "2026/12/29")
# elsewhere in the sheet, so force the cell to Text before assigning the value to
# avoid Excel auto-converting the "YYYY/MM/DD" string into a date serial, then
# clear the format again so no stray number-format style is left behind.
$ws.Range("A671").NumberFormat = "@"
$ws.Range("A671").Value = "2026/01/21"
$ws.Range("A671").ClearFormats()

$ws.Range("B671").Value = "水"
$ws.Range("C671").Value = 19
$ws.Range("D671").Value = 27
